$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column B (the "Run" boolean flag) to TRUE for every scenario row (3-81).
# Row 46 is already TRUE in the source file; re-setting it is harmless.
for ($r = 3; $r -le 81; $r++) {
    $ws.Range("B$r").Value = $true
}

# Update the saved view/cursor state: the user scrolled the frozen pane down
# and left the selection on C61 (bottom-right pane), while the freeze itself
# (columns/rows 1:2 frozen) is unchanged.
$ws.Range("C61").Select()
